$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.321.40'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '2.987.62'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.06'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.50'
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '2.987.25'
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.05'
$ws.Range("E11").Value = '  +7.11%  '
$ws.Range("E12").Value = '  +3.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").Value = '3.479.05'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").Value = '61.307.37'
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("D19").Value = '2.982.87'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.82'
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.07'
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.684'
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.31'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.54'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.66'
$ws.Range("E25").Value = '  +5.72%  '
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.90'
$ws.Range("E27").Value = '  -2.95%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.21'
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("E35").Value = '  +3.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.42'
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.14'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.06'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("E41").Value = '  +10.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '397.36'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.72'
$ws.Range("E44").Value = '  +4.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0351'
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  -2.68%  '
$ws.Range("D47").Value = '2.689.44'
$ws.Range("E47").Value = '  -2.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.79'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("E51").Value = '  +1.17%  '
